# Generate Report for Handoff
# Adds two new handoff rows (e3a9b16a-... and fd0329ee-...) to the
# Overview / zh-cn / de-de worksheets of the localization-status workbook.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276   # OLE color for RGB(0x64,0x95,0xED) == FF6495ED

function Set-HyperlinkFont($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

function Set-DateFormat($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Cells.Item(4, 1).Value = "e3a9b16a-1b58-4003-9522-9be81fb79950.md"
$ws1.Cells.Item(4, 2).Value = "Ready for handoff"
$ws1.Cells.Item(4, 3).Value = "Ready for handoff"
$ws1.Cells.Item(4, 4).Value = "2016-03-23 18:42:44"

$ws1.Cells.Item(5, 1).Value = "fd0329ee-f51a-4870-bb82-0f99a9481b66.md"
$ws1.Cells.Item(5, 2).Value = "Ready for handoff"
$ws1.Cells.Item(5, 3).Value = "Ready for handoff"
$ws1.Cells.Item(5, 4).Value = "2016-03-23 18:42:44"

$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e3a9b16a/e2e/e3a9b16a-1b58-4003-9522-9be81fb79950.md", "", "", "e3a9b16a-1b58-4003-9522-9be81fb79950.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/fd0329ee/e2e/fd0329ee-f51a-4870-bb82-0f99a9481b66.md", "", "", "fd0329ee-f51a-4870-bb82-0f99a9481b66.md") | Out-Null

Set-HyperlinkFont $ws1.Range("A4")
Set-HyperlinkFont $ws1.Range("A5")
Set-DateFormat $ws1.Range("D4")
Set-DateFormat $ws1.Range("D5")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Cells.Item(4, 1).Value = "e3a9b16a-1b58-4003-9522-9be81fb79950.md"
$ws2.Cells.Item(4, 2).Value = ".md"
$ws2.Cells.Item(4, 3).Value = "Ready for handoff"
$ws2.Cells.Item(4, 4).Value = "e3a9b16a-1b58-4003-9522-9be81fb79950.8ea1b1c973f2d86da1c4c742ae52735de00dfee0.zh-cn.xlf"
$ws2.Cells.Item(4, 5).Value = "2016-03-23 18:42:40"
$ws2.Cells.Item(4, 8).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(4, 10).Value = "Include"

$ws2.Cells.Item(5, 1).Value = "fd0329ee-f51a-4870-bb82-0f99a9481b66.md"
$ws2.Cells.Item(5, 2).Value = ".md"
$ws2.Cells.Item(5, 3).Value = "Ready for handoff"
$ws2.Cells.Item(5, 4).Value = "fd0329ee-f51a-4870-bb82-0f99a9481b66.1da723cd92c9571137e1079134d51d53b326f9ce.zh-cn.xlf"
$ws2.Cells.Item(5, 5).Value = "2016-03-23 18:42:40"
$ws2.Cells.Item(5, 8).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(5, 10).Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e3a9b16a/e2e/e3a9b16a-1b58-4003-9522-9be81fb79950.md", "", "", "e3a9b16a-1b58-4003-9522-9be81fb79950.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e3a9b16a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e3a9b16a-1b58-4003-9522-9be81fb79950.8ea1b1c973f2d86da1c4c742ae52735de00dfee0.zh-cn.xlf", "", "", "e3a9b16a-1b58-4003-9522-9be81fb79950.8ea1b1c973f2d86da1c4c742ae52735de00dfee0.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/fd0329ee/e2e/fd0329ee-f51a-4870-bb82-0f99a9481b66.md", "", "", "fd0329ee-f51a-4870-bb82-0f99a9481b66.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fd0329ee/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/fd0329ee-f51a-4870-bb82-0f99a9481b66.1da723cd92c9571137e1079134d51d53b326f9ce.zh-cn.xlf", "", "", "fd0329ee-f51a-4870-bb82-0f99a9481b66.1da723cd92c9571137e1079134d51d53b326f9ce.zh-cn.xlf") | Out-Null

Set-HyperlinkFont $ws2.Range("A4")
Set-HyperlinkFont $ws2.Range("D4")
Set-HyperlinkFont $ws2.Range("A5")
Set-HyperlinkFont $ws2.Range("D5")
Set-DateFormat $ws2.Range("E4")
Set-DateFormat $ws2.Range("H4")
Set-DateFormat $ws2.Range("E5")
Set-DateFormat $ws2.Range("H5")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Cells.Item(4, 1).Value = "e3a9b16a-1b58-4003-9522-9be81fb79950.md"
$ws3.Cells.Item(4, 2).Value = ".md"
$ws3.Cells.Item(4, 3).Value = "Ready for handoff"
$ws3.Cells.Item(4, 4).Value = "e3a9b16a-1b58-4003-9522-9be81fb79950.8ea1b1c973f2d86da1c4c742ae52735de00dfee0.de-de.xlf"
$ws3.Cells.Item(4, 5).Value = "2016-03-23 18:42:44"
$ws3.Cells.Item(4, 8).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(4, 10).Value = "Include"

$ws3.Cells.Item(5, 1).Value = "fd0329ee-f51a-4870-bb82-0f99a9481b66.md"
$ws3.Cells.Item(5, 2).Value = ".md"
$ws3.Cells.Item(5, 3).Value = "Ready for handoff"
$ws3.Cells.Item(5, 4).Value = "fd0329ee-f51a-4870-bb82-0f99a9481b66.1da723cd92c9571137e1079134d51d53b326f9ce.de-de.xlf"
$ws3.Cells.Item(5, 5).Value = "2016-03-23 18:42:44"
$ws3.Cells.Item(5, 8).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(5, 10).Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e3a9b16a/e2e/e3a9b16a-1b58-4003-9522-9be81fb79950.md", "", "", "e3a9b16a-1b58-4003-9522-9be81fb79950.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e3a9b16a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e3a9b16a-1b58-4003-9522-9be81fb79950.8ea1b1c973f2d86da1c4c742ae52735de00dfee0.de-de.xlf", "", "", "e3a9b16a-1b58-4003-9522-9be81fb79950.8ea1b1c973f2d86da1c4c742ae52735de00dfee0.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/fd0329ee/e2e/fd0329ee-f51a-4870-bb82-0f99a9481b66.md", "", "", "fd0329ee-f51a-4870-bb82-0f99a9481b66.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fd0329ee/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/fd0329ee-f51a-4870-bb82-0f99a9481b66.1da723cd92c9571137e1079134d51d53b326f9ce.de-de.xlf", "", "", "fd0329ee-f51a-4870-bb82-0f99a9481b66.1da723cd92c9571137e1079134d51d53b326f9ce.de-de.xlf") | Out-Null

Set-HyperlinkFont $ws3.Range("A4")
Set-HyperlinkFont $ws3.Range("D4")
Set-HyperlinkFont $ws3.Range("A5")
Set-HyperlinkFont $ws3.Range("D5")
Set-DateFormat $ws3.Range("E4")
Set-DateFormat $ws3.Range("H4")
Set-DateFormat $ws3.Range("E5")
Set-DateFormat $ws3.Range("H5")

Write-Host "Report generated for handoff."
